$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the canonical 2018 dates for the "expected end" / "actual start" /
# "actual end" columns (G, H, I). Rows 4-10 had accidentally been filled with a
# staircase of future years (2019-2025) instead of repeating the 2018 date -
# reset them back to match row 3.
$ws.Range("G4:G10").Value = "30/9/2018"
$ws.Range("H4:H10").Value = "19/9/2018"
$ws.Range("I4:I10").Value = "30/9/2018"

# Leave the cursor where the user finished editing.
$ws.Range("H7").Select() | Out-Null
